# Little Library (Albuquerque) data export
# "changed PostgreSQL file slightly to fit prior push to Library Spec dropdown"
#
# - Rename the worksheet to match the exported table name.
# - Drop a couple of stray "library_specialization" values that no longer
#   belong (a duplicated library-type value, and the removed "Science
#   Fiction" category).
# - Rename "Young Adults" to "Young Adult" to match the updated dropdown.
# - Restore the cursor/selection position that was active when the file
#   was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from the generic default name to the table name.
$ws.Name = "library-abq-all"

# Row 2: library_specialization erroneously held the library type value
# ("Little Library"); clear it so only library_type (column J) has it.
$ws.Range("I2").Clear()

# Row 22: the "Science Fiction" specialization option was removed from the
# Library Spec dropdown, so clear the stale reference.
$ws.Range("I22").Clear()

# Row 27: "Young Adults" was renamed to "Young Adult" in the dropdown.
$ws.Range("I27").Value = "Young Adult"

# Restore the saved view/selection (user had scrolled back up and selected
# E28 before saving).
$ws.Range("E28").Select()
